# CS5044-Report.docx edit:
# "Added the number of observations of the data"
#
# The draft text reads:
#   "... which are ???? from the original ???? observations"
# and the author filled in the actual observation counts:
#   "... which are 3339 from the original 4700 observations"
#
# wdFindContinue   = 1
# wdReplaceOne     = 1
# wdReplaceAll     = 2

$d = $word.ActiveDocument

# First "????" -> size of the 5-league subsample (3339 observations)
$range1 = $d.Content
$range1.Find.Execute("????", $false, $false, $false, $false, $false, $true, 1, $false, "3339", 1)

# Second "????" -> size of the original full dataset (4700 observations)
$range2 = $d.Content
$range2.Find.Execute("????", $false, $false, $false, $false, $false, $true, 1, $false, "4700", 1)
